$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "technical data layer" columns (C: raw value, D: difference vs column B)
$ws.Range("C2").Value = 0.061
$ws.Range("D2").Formula = "=C2-B2"

$ws.Range("C3").Value = 0.055
$ws.Range("D3").Formula = "=C3-B3"

# Matches the cursor position left behind in the authored workbook
$ws.Range("D4").Select()
